{"js": "// Each entry is [oldText, newText]. The document contains exactly one\n// occurrence of each oldText (the date line plus 25 division problems),\n// so a literal, case-sensitive search-and-replace is safe and unambiguous.\nconst replacements = [\n  [\"2024-12-24 Tuesday\", \"2024-12-25 Wednesday\"],\n  [\"736\u00f77=105, 1\", \"176\u00f77=25, 1\"],\n  [\"538\u00f79=59, 7\", \"587\u00f79=65, 2\"],\n  [\"350\u00f72=175, 0\", \"766\u00f75=153, 1\"],\n  [\"886\u00f75=177, 1\", \"428\u00f75=85, 3\"],\n  [\"319\u00f77=45, 4\", \"987\u00f78=123, 3\"],\n  [\"519\u00f74=129, 3\", \"640\u00f78=80, 0\"],\n  [\"445\u00f75=89, 0\", \"935\u00f74=233, 3\"],\n  [\"163\u00f75=32, 3\", \"466\u00f79=51, 7\"],\n  [\"607\u00f75=121, 2\", \"239\u00f77=34, 1\"],\n  [\"951\u00f73=317, 0\", \"585\u00f76=97, 3\"],\n  [\"976\u00f79=108, 4\", \"556\u00f74=139, 0\"],\n  [\"382\u00f76=63, 4\", \"807\u00f72=403, 1\"],\n  [\"722\u00f77=103, 1\", \"462\u00f73=154, 0\"],\n  [\"656\u00f78=82, 0\", \"962\u00f76=160, 2\"],\n  [\"261\u00f74=65, 1\", \"504\u00f76=84, 0\"],\n  [\"884\u00f73=294, 2\", \"455\u00f78=56, 7\"],\n  [\"865\u00f74=216, 1\", \"526\u00f74=131, 2\"],\n  [\"346\u00f77=49, 3\", \"868\u00f78=108, 4\"],\n  [\"997\u00f72=498, 1\", \"788\u00f72=394, 0\"],\n  [\"841\u00f75=168, 1\", \"453\u00f75=90, 3\"],\n  [\"289\u00f77=41, 2\", \"874\u00f73=291, 1\"],\n  [\"314\u00f77=44, 6\", \"983\u00f77=140, 3\"],\n  [\"274\u00f73=91, 1\", \"603\u00f76=100, 3\"],\n  [\"515\u00f74=128, 3\", \"262\u00f77=37, 3\"],\n  [\"963\u00f72=481, 1\", \"600\u00f78=75, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Each pair is (oldText, newText). The document contains exactly one\n# occurrence of each oldText (the date line plus 25 division problems),\n# so a literal Find/Replace (MatchCase, no wildcards) is safe and unambiguous.\n$d = $word.ActiveDocument\n$pairs = @(\n  @(\"2024-12-24 Tuesday\", \"2024-12-25 Wednesday\"),\n  @(\"736\u00f77=105, 1\", \"176\u00f77=25, 1\"),\n  @(\"538\u00f79=59, 7\", \"587\u00f79=65, 2\"),\n  @(\"350\u00f72=175, 0\", \"766\u00f75=153, 1\"),\n  @(\"886\u00f75=177, 1\", \"428\u00f75=85, 3\"),\n  @(\"319\u00f77=45, 4\", \"987\u00f78=123, 3\"),\n  @(\"519\u00f74=129, 3\", \"640\u00f78=80, 0\"),\n  @(\"445\u00f75=89, 0\", \"935\u00f74=233, 3\"),\n  @(\"163\u00f75=32, 3\", \"466\u00f79=51, 7\"),\n  @(\"607\u00f75=121, 2\", \"239\u00f77=34, 1\"),\n  @(\"951\u00f73=317, 0\", \"585\u00f76=97, 3\"),\n  @(\"976\u00f79=108, 4\", \"556\u00f74=139, 0\"),\n  @(\"382\u00f76=63, 4\", \"807\u00f72=403, 1\"),\n  @(\"722\u00f77=103, 1\", \"462\u00f73=154, 0\"),\n  @(\"656\u00f78=82, 0\", \"962\u00f76=160, 2\"),\n  @(\"261\u00f74=65, 1\", \"504\u00f76=84, 0\"),\n  @(\"884\u00f73=294, 2\", \"455\u00f78=56, 7\"),\n  @(\"865\u00f74=216, 1\", \"526\u00f74=131, 2\"),\n  @(\"346\u00f77=49, 3\", \"868\u00f78=108, 4\"),\n  @(\"997\u00f72=498, 1\", \"788\u00f72=394, 0\"),\n  @(\"841\u00f75=168, 1\", \"453\u00f75=90, 3\"),\n  @(\"289\u00f77=41, 2\", \"874\u00f73=291, 1\"),\n  @(\"314\u00f77=44, 6\", \"983\u00f77=140, 3\"),\n  @(\"274\u00f73=91, 1\", \"603\u00f76=100, 3\"),\n  @(\"515\u00f74=128, 3\", \"262\u00f77=37, 3\"),\n  @(\"963\u00f72=481, 1\", \"600\u00f78=75, 0\")\n)\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n  #          MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n  # MatchCase=$true, Wrap=wdFindContinue(1), Replace=wdReplaceAll(2)\n  $result = $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n  if (-not $result) {\n    throw \"Replacement failed for: $old\"\n  }\n}"}
